$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7822278738021851
$ws.Range("B1").Value = 1.247182250022888
$ws.Range("C1").Value = 3.148746252059937
$ws.Range("D1").Value = 3.121608018875122
$ws.Range("E1").Value = 1.722864985466003
